$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New roster data (player, position, team) for rows 2-19
$data = @(
    @("Anfernee Simons", "PG,SG", "Portland Trail Blazers"),
    @("Anthony Edwards", "SG,SF", "Minnesota Timberwolves"),
    @("Amen Thompson", "SG,SF,PF", "Houston Rockets"),
    @("Jerami Grant", "SF,PF", "Portland Trail Blazers"),
    @("Jaren Jackson Jr.", "PF,C", "Memphis Grizzlies"),
    @("Bradley Beal", "PG,SG,SF", "Phoenix Suns"),
    @("Jayson Tatum", "SF,PF", "Boston Celtics"),
    @("Zion Williamson", "PF,C", "New Orleans Pelicans"),
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("Zach Edey", "C", "Memphis Grizzlies"),
    @("Giannis Antetokounmpo", "PF,C", "Milwaukee Bucks"),
    @("Fred VanVleet", "PG", "Houston Rockets"),
    @("James Harden", "PG,SG", "LA Clippers"),
    @("RJ Barrett", "SG,SF,PF", "Toronto Raptors"),
    @("Ivica Zubac", "C", "LA Clippers"),
    @("Jonathan Kuminga", "SF,PF", "Golden State Warriors"),
    @("Paul George", "SG,SF,PF", "Philadelphia 76ers"),
    @("P.J. Washington", "SF,PF", "Dallas Mavericks")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
